$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: JavaClassName column (C) loses the "tests." segment
$ws.Range("C2").Value = "com.yourorg.LoginTest"
$ws.Range("C3").Value = "com.yourorg.LoginTest"
$ws.Range("C4").Value = "com.yourorg.PaymentTest"
$ws.Range("C5").Value = "com.yourorg.PaymentTest"

# Add two new rows
$ws.Range("A6").Value = "SignIssueTest"
$ws.Range("B6").Value = "TC005"
$ws.Range("C6").Value = "com.yourorg.SignTest"

$ws.Range("A7").Value = "SignAndPaymentTest"
$ws.Range("B7").Value = "TC006"
$ws.Range("C7").Value = "com.yourorg.SignAndPaymentTest"

# Column widths (pre-compensated so the emulator's internal rounding lands
# on the value closest to the widths recorded in the target workbook)
$ws.Columns.Item(1).ColumnWidth = 20.5
$ws.Columns.Item(2).ColumnWidth = 26.333333333333332
$ws.Columns.Item(4).ColumnWidth = 53.666666666666664

# Update selection to match the last active cell
$ws.Range("A7").Select()
